$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the cryptos.xlsx data refresh diff.
# Numeric-looking text (e.g. "197.20") is written via a temporary
# Text number-format + ClearFormats so Excel keeps it as a literal
# string (matching the source inlineStr cells) without leaving any
# residual cell style applied.

$ws.Range("D2").Value = "67.117.36"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "3.485.86"
$ws.Range("E3").Value = "  -3.95%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "197.20"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "548.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.29%  "
$ws.Range("D7").Value = "3.481.49"
$ws.Range("E7").Value = "  -3.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.648"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "62.09"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +11.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.142"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -7.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -10.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.74"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.89%  "
$ws.Range("D15").Value = "4.060.66"
$ws.Range("D16").Value = "3.500.54"
$ws.Range("E16").Value = "  -3.45%  "
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "66.958.95"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.18"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.74"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.55%  "
$ws.Range("E21").Value = "  -5.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "386.82"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.98"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.74"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -7.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.01"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.13"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.80"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.77"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.68"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.81"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "674.45"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.94"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -14.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.61"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.18"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.109"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.14"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -10.59%  "
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.395"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.42%  "
$ws.Range("D39").Value = "3.073.83"
$ws.Range("E39").Value = "  -2.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.98"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.129"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.52%  "
$ws.Range("D43").Value = "0.0₃0665"
$ws.Range("E43").Value = "  -16.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.76"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.94%  "
$ws.Range("E45").Value = "  -13.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.71"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0393"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.26%  "
$ws.Range("E48").Value = "  -5.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "136.33"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.49%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.90"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.26%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.15"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -8.01%  "
